# Populate the header row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Work Updates"

# Widen column B to fit the "Work Updates" header
$ws.Columns.Item(2).ColumnWidth = 18.92

# Leave the selection where the author left it
$ws.Range("C5").Select()
